$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the daily conversion note text (A1 on Hoja1)
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.32 = 8719.29 pesos`n✅ 8719.29 pesos = 2.31 = 950.82 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# Update the "tasas" rate values
$ws2.Range("N10").Value = 430.776
$ws2.Range("O10").Value = 3756.06
$ws2.Range("N12").Value = 3780
$ws2.Range("O12").Value = 412.2
